# Implement Common data provider, soft assertions, parameterization
#
# Adds a second worksheet ("OpenAccountTest") after the existing
# "AddCustomerTest" sheet, populates it with customer/currency test data,
# and appends two more parameterized rows to the "AddCustomerTest" sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- New sheet: OpenAccountTest, inserted right after AddCustomerTest ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "OpenAccountTest"

$ws2.Range("A1").Value = "customer"
$ws2.Range("B1").Value = "currency"
$ws2.Range("A2").Value = "Dan Car"
$ws2.Range("B2").Value = "Rupee"

[void]$ws2.Range("D8").Select()

# --- More parameterized rows on AddCustomerTest ---
$ws1.Range("A3").Value = "Sam"
$ws1.Range("B3").Value = "Tar"
$ws1.Range("A4").Value = "Pam"
$ws1.Range("B4").Value = "Sen"
$ws1.Range("C3").Value = "syd2143"
$ws1.Range("C4").Value = "syd2146"
$ws1.Range("D3").Value = "Customer added successfully"
$ws1.Range("D4").Value = "Customer added successfully"

# Leave the original sheet active with the cursor on the last data row.
[void]$ws1.Activate()
[void]$ws1.Range("C4").Select()
